# Backlog Projet Création site web RIL10.xlsx
# Feature: "suppression groupée" (bulk delete) user story moved from
# TODO -> DONE and assigned to JB, along with a few other already-DOING
# stories being marked DONE (assigned to JB) as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille 1")

# Row 2 : Story "ajouter une intervention" - was TODO(x)/DOING(JB) -> now DOING(JB)/DONE(x)
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = "x"

# Row 3 : Story "éditer une intervention" - was DONE(x) only -> now DOING(JB)/DONE(x)
$ws.Range("E3").Value = "JB"

# Row 5 : Story "visualiser les détails d'une intervention" - was TODO(x) -> now DOING(JB)/DONE(x)
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "JB"
$ws.Range("F5").Value = "x"

# Row 11 : Story "supprimer une intervention" (groupée) - was TODO(x) -> now DOING(JB)/DONE(x)
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = "JB"
$ws.Range("F11").Value = "x"

# Row 12 : Story "limiter le nombre d'interventions affichées" - was TODO(x) -> now DOING(JB)/DONE(x)
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "JB"
$ws.Range("F12").Value = "x"

# Move the active selection to E14 (last user action before save)
$ws.Range("E14").Select()
